$d = $word.ActiveDocument

# The document has two sections worth of headers/footers, each holding the
# same pair of logo pictures (BTec_Logo-Orange in the headers, the Pearson
# logo in the footers). Each picture's non-visual drawing properties carry a
# cosmetic "name" (wp:docPr / pic:cNvPr @name) that needs to be swapped:
#   - headers:  image1.jpg <-> image2.jpg   (BTec_Logo-Orange)
#   - footers:  image2.png <-> image1.png   (Pearson logo)
#
# InlineShape.Name has to be set through the Selection (rather than the
# InlineShape reference straight off Range.InlineShapes) so the edit resolves
# reliably for shapes living in footer stories as well as header stories.

function Rename-InlineLogo {
    param($shapeRange, [string]$newName)

    $shapeRange.Select()
    $word.Selection.InlineShapes(1).Name = $newName
}

# Headers -- BTec_Logo-Orange: image1.jpg -> image2.jpg
Rename-InlineLogo $d.Sections(1).Headers(2).Range.InlineShapes(1) "image2.jpg"
Rename-InlineLogo $d.Sections(1).Headers(1).Range.InlineShapes(1) "image2.jpg"

# Footers -- Pearson logo: image2.png -> image1.png
Rename-InlineLogo $d.Sections(1).Footers(2).Range.InlineShapes(1) "image1.png"
Rename-InlineLogo $d.Sections(1).Footers(1).Range.InlineShapes(1) "image1.png"
